# Update "想去人数" (interested-people count) figures in the
# "展览" (sheet1), "演出" (sheet2) and "全部类型" (sheet4) worksheets
# to match the newly generated data snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws2 = $wb.Worksheets.Item("演出")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- 展览 sheet updates ---
$ws1.Range("F3").Value = 163
$ws1.Range("F5").Value = 47
$ws1.Range("F6").Value = 2734
$ws1.Range("F8").Value = 1618
$ws1.Range("F9").Value = 7408
$ws1.Range("F11").Value = 7590
$ws1.Range("F15").Value = 6087
$ws1.Range("F16").Value = 3241
$ws1.Range("F17").Value = 3615
$ws1.Range("F19").Value = 5
$ws1.Range("F21").Value = 28
$ws1.Range("F22").Value = 441
$ws1.Range("F26").Value = 2102
$ws1.Range("F31").Value = 1071
$ws1.Range("F32").Value = 60
$ws1.Range("F33").Value = 13
$ws1.Range("F34").Value = 2591
$ws1.Range("F35").Value = 1438
$ws1.Range("F37").Value = 10
$ws1.Range("F38").Value = 18
$ws1.Range("F39").Value = 3205
$ws1.Range("F40").Value = 147
$ws1.Range("F43").Value = 892
$ws1.Range("F44").Value = 474
$ws1.Range("F45").Value = 1254
$ws1.Range("F48").Value = 581

# --- 演出 sheet updates ---
$ws2.Range("F4").Value = 54

# --- 全部类型 sheet updates ---
$ws4.Range("F5").Value = 163
$ws4.Range("F6").Value = 54
$ws4.Range("F7").Value = 47
$ws4.Range("F9").Value = 2734
$ws4.Range("F10").Value = 1618
$ws4.Range("F13").Value = 7408
$ws4.Range("F14").Value = 7590
$ws4.Range("F17").Value = 6087
$ws4.Range("F18").Value = 3241
$ws4.Range("F19").Value = 3615
$ws4.Range("F21").Value = 5
$ws4.Range("F22").Value = 28
$ws4.Range("F23").Value = 441
$ws4.Range("F29").Value = 2102
$ws4.Range("F37").Value = 13
$ws4.Range("F38").Value = 2591
$ws4.Range("F39").Value = 1438
$ws4.Range("F41").Value = 10
$ws4.Range("F43").Value = 3205
$ws4.Range("F45").Value = 892
$ws4.Range("F46").Value = 474
$ws4.Range("F47").Value = 1254
$ws4.Range("F49").Value = 581
